# Update Korean New Deal workbook:
#  - Insert a new "EDA" worksheet (a copy of "finance") right after "finance"
#  - Shorten four of its category labels (column D) so they read as the
#    new, more concise descriptions
#  - Make "EDA" the active/selected sheet (finance keeps its old selection,
#    but is no longer the tab shown on open)

$wb = $excel.ActiveWorkbook

$finance = $wb.Worksheets.Item("finance")

# Duplicate "finance" and drop the copy immediately after it.
$finance.Copy($null, $finance)
$eda = $wb.ActiveSheet
$eda.Name = "EDA"

# Shorten the category descriptions in column D (order matters so the new
# shared-string table is built in the same sequence as the source edit).
$eda.Range("D2").Value = "데이터 구축ㆍ개방ㆍ활용"
$eda.Range("D14").Value = "공공시설 제로에너지화"
$eda.Range("D6").Value = "초중고에 디지털 기반 교육 인프라 조성"
$eda.Range("D7").Value = "대학ㆍ직업훈련기관 온라인 교육 강화"

# Restore the working selection for the new sheet and make it the active tab.
$null = $eda.Range("D30").Select()
